# "Generate Report for Archive"
#
# 1. The localization status moved on from "Ready for handoff" to
#    "In Translation" everywhere it is shown (Overview!E2:F3,
#    zh-cn!C2:C3, de-de!C2:C3 all reference the same shared string).
# 2. The Status-column widths on the Overview sheet (E:F) and on each
#    per-locale sheet (C) were re-sized narrower as part of regenerating
#    the report.

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text everywhere it appears -----------------
foreach ($ws in $wb.Worksheets) {
    $ws.UsedRange.Replace("Ready for handoff", "In Translation") | Out-Null
}

# --- 2. Narrow the Status column(s) -----------------------------------
# Target stored column width (OOXML <col width=.../>) is 13.4101845877511
# characters. ColumnWidth is also specified in characters, so set it to
# the ColumnWidth that maps to the stored width closest to that target.
$newStatusWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newStatusWidth   # column E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = $newStatusWidth   # column F (de-de)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newStatusWidth        # column C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newStatusWidth        # column C (Status)
